$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("python")

# Update B64: "leave" -> "Assignment(13)"
$ws.Range("B64").Value = "Assignment(13)"

# Clear B65 and B66 (cells removed in the diff)
$ws.Range("B65").ClearContents()
$ws.Range("B66").ClearContents()

# Update the view: scroll position and active selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 52
$ws.Range("B64").Select()
